$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be written as literal text, matching the
# original inline-string cells (values like "519.01" or "0.998" would
# otherwise be auto-coerced to numbers by Excel's smart-entry parsing).
$priceCells = @("D2", "D3", "D5", "D9", "D11", "D12", "D14", "D15", "D16", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D32", "D33", "D37", "D40", "D42", "D43", "D44", "D48", "D49", "D50", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "56.319.88"
$ws.Range("E2").Value = "  +3.81%  "

$ws.Range("D3").Value = "2.316.68"
$ws.Range("E3").Value = "  +2.36%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "519.01"
$ws.Range("E5").Value = "  +4.71%  "

$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("E8").Value = "  +2.15%  "

$ws.Range("D9").Value = "2.341.65"
$ws.Range("E9").Value = "  +3.22%  "

$ws.Range("E10").Value = "  +8.28%  "

$ws.Range("D11").Value = "0.155"
$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").Value = "5.18"
$ws.Range("E12").Value = "  +7.29%  "

$ws.Range("E13").Value = "  +2.15%  "

$ws.Range("D14").Value = "24.10"
$ws.Range("E14").Value = "  +5.14%  "

$ws.Range("D15").Value = "2.731.86"
$ws.Range("E15").Value = "  +2.53%  "

$ws.Range("D16").Value = "56.417.38"
$ws.Range("E16").Value = "  +4.06%  "

$ws.Range("E17").Value = "  +4.86%  "

$ws.Range("D18").Value = "2.353.47"
$ws.Range("E18").Value = "  +3.80%  "

$ws.Range("E19").Value = "  +3.26%  "

$ws.Range("E20").Value = "  +3.98%  "

$ws.Range("D21").Value = "321.20"
$ws.Range("E21").Value = "  +6.05%  "

$ws.Range("D22").Value = "6.64"
$ws.Range("E22").Value = "  +4.96%  "

$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "60.82"
$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("D25").Value = "0.993"
$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("D26").Value = "0.160"
$ws.Range("E26").Value = "  +7.16%  "

$ws.Range("D27").Value = "7.68"
$ws.Range("E27").Value = "  +4.97%  "

$ws.Range("D28").Value = "171.52"
$ws.Range("E28").Value = "  +0.62%  "

$ws.Range("E29").Value = "  +12.33%  "

$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  +6.97%  "

$ws.Range("E31").Value = "  +5.92%  "

$ws.Range("D32").Value = "6.28"
$ws.Range("E32").Value = "  +5.27%  "

$ws.Range("D33").Value = "18.43"
$ws.Range("E33").Value = "  +3.69%  "

$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("E35").Value = "  -0.57%  "

$ws.Range("E36").Value = "  +5.86%  "

$ws.Range("D37").Value = "0.936"
$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("E39").Value = "  +9.12%  "

$ws.Range("D40").Value = "37.49"
$ws.Range("E40").Value = "  +4.51%  "

$ws.Range("E41").Value = "  +2.34%  "

$ws.Range("D42").Value = "140.99"
$ws.Range("E42").Value = "  +13.05%  "

$ws.Range("D43").Value = "3.60"
$ws.Range("E43").Value = "  +6.94%  "

$ws.Range("D44").Value = "280.97"
$ws.Range("E44").Value = "  +16.38%  "

$ws.Range("E45").Value = "  +6.27%  "

$ws.Range("E46").Value = "  +3.65%  "

$ws.Range("E47").Value = "  +3.97%  "

$ws.Range("D48").Value = "0.559"
$ws.Range("E48").Value = "  +2.58%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0216"
$ws.Range("E49").Value = "  +6.03%  "

$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D50").Value = "0.382"
$ws.Range("E50").Value = "  +2.34%  "

$ws.Range("D51").Value = "16.98"
$ws.Range("E51").Value = "  +5.47%  "

# Restore the Price column's number format back to General, same as
# every other untouched cell in the sheet (the text value set above is
# preserved -- only the display/number-format reverts).
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "General"
}